$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 18; rows 18-25 shift down to 19-26,
# and the inserted row inherits formatting (incl. the date style) from
# the row above it, same as Excel's native "Insert Row" behavior.
$ws.Rows(18).Insert()

# Populate the new row 18 with its own data (a new weekly price record).
$ws.Range("A18").Value = 1
$ws.Range("B18").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C18").Value = "Arica y Parinacota"
$ws.Range("D18").Value = 44566
$ws.Range("E18").Value = 15
$ws.Range("F18").Value = 100112028
$ws.Range("G18").Value = "Sandia"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 1000
$ws.Range("K18").Value = 300
$ws.Range("L18").Value = 320
$ws.Range("M18").Value = 310
$ws.Range("N18").Value = "$/kilo (volumen en unidades)"
$ws.Range("O18").Value = "Perú"
$ws.Range("P18").Value = 310
$ws.Range("Q18").Value = 1
$ws.Range("R18").Value = "Hortaliza"
